$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.129.96"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "2.472.62"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").Value = "2.472.03"
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "2.924.13"
$ws.Range("E16").Value = "  +2.59%  "

$ws.Range("D17").Value = "63.135.75"
$ws.Range("E17").Value = "  +1.66%  "

$ws.Range("D18").Value = "2.477.66"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "668.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.84%  "

$ws.Range("D28").Value = "0.0₃0994"
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").Value = "2.594.91"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1,128.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.21%  "

$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("E35").Value = "  +4.60%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("E40").Value = "  +0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").Value = "0.0₆0306"
$ws.Range("E45").Value = "  +11.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +27.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.607"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0515"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "

